$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Experiments")

# ---- Row 36: new section header (NYU Depth v2 Labeled Test Set (654 Entries)) ----
$ws.Cells.Item(36, 1).Value = "NYU Depth v2 Labled Test Set (654 Entries)"
$ws.Cells.Item(36, 2).NumberFormat = "0.000"
$ws.Cells.Item(36, 2).Value = "delta1"
$ws.Cells.Item(36, 3).NumberFormat = "0.000"
$ws.Cells.Item(36, 3).Value = "delta2"
$ws.Cells.Item(36, 4).NumberFormat = "0.000"
$ws.Cells.Item(36, 4).Value = "delta3"
$ws.Cells.Item(36, 5).NumberFormat = "0.000"
$ws.Cells.Item(36, 5).Value = "mse"
$ws.Cells.Item(36, 6).NumberFormat = "0.000"
$ws.Cells.Item(36, 6).Value = "rmse"
$ws.Cells.Item(36, 7).NumberFormat = "0.000"
$ws.Cells.Item(36, 7).Value = "rel_abs_dif"
$ws.Cells.Item(36, 8).NumberFormat = "0.000"
$ws.Cells.Item(36, 8).Value = "rel_sqr_diff"
$ws.Cells.Item(36, 9).NumberFormat = "0.000"
$ws.Cells.Item(36, 9).Value = "log10"
$ws.Cells.Item(36, 10).NumberFormat = "0.000"
$ws.Cells.Item(36, 10).Value = "log_rmse"

# ---- Row 37: DORN_nohints results on the labeled test set ----
$ws.Cells.Item(37, 1).Value = "DORN_nohints"
$ws.Cells.Item(37, 2).NumberFormat = "0.000"
$ws.Cells.Item(37, 2).Value = 0.83929108669991404
$ws.Cells.Item(37, 3).NumberFormat = "0.000"
$ws.Cells.Item(37, 3).Value = 0.95858464605035199
$ws.Cells.Item(37, 4).NumberFormat = "0.000"
$ws.Cells.Item(37, 4).Value = 0.98564354107743601
$ws.Cells.Item(37, 5).NumberFormat = "0.000"
$ws.Cells.Item(37, 5).Value = 0.26140680609123301
$ws.Cells.Item(37, 6).NumberFormat = "0.000"
$ws.Cells.Item(37, 6).Formula = "=SQRT(E37)"
$ws.Cells.Item(37, 7).NumberFormat = "0.000"
$ws.Cells.Item(37, 7).Value = 0.12950759483657001
$ws.Cells.Item(37, 8).NumberFormat = "0.000"
$ws.Cells.Item(37, 8).Value = 0.0875682236801333
$ws.Cells.Item(37, 9).NumberFormat = "0.000"
$ws.Cells.Item(37, 9).Value = 0.0585035172251433
$ws.Cells.Item(37, 10).NumberFormat = "0.000"
$ws.Cells.Item(37, 10).Value = 0.17285743104927601

# ---- Row 38: DenseDepth_nohints results on the labeled test set ----
$ws.Cells.Item(38, 1).Value = "DenseDepth_nohints"
$ws.Cells.Item(38, 2).NumberFormat = "0.000"
$ws.Cells.Item(38, 2).Value = 0.85618639232578697
$ws.Cells.Item(38, 3).NumberFormat = "0.000"
$ws.Cells.Item(38, 3).Value = 0.97836328187102095
$ws.Cells.Item(38, 4).NumberFormat = "0.000"
$ws.Cells.Item(38, 4).Value = 0.99565340722963402
$ws.Cells.Item(38, 5).NumberFormat = "0.000"
$ws.Cells.Item(38, 5).Value = 0.21353111323962001
$ws.Cells.Item(38, 6).NumberFormat = "0.000"
$ws.Cells.Item(38, 6).Formula = "=SQRT(E38)"
$ws.Cells.Item(38, 7).NumberFormat = "0.000"
$ws.Cells.Item(38, 7).Value = 0.11950074903631699
$ws.Cells.Item(38, 8).NumberFormat = "0.000"
$ws.Cells.Item(38, 8).Value = 0.0687247217789462
$ws.Cells.Item(38, 9).NumberFormat = "0.000"
$ws.Cells.Item(38, 9).Value = 0.0513834079050701
$ws.Cells.Item(38, 10).NumberFormat = "0.000"
$ws.Cells.Item(38, 10).Value = 0.151374158361015

# ---- Notes column (L): written last, row 38 first then row 37 ----
$ws.Cells.Item(38, 12).WrapText = $true
$ws.Cells.Item(38, 12).Value = "Uses rawdepth for evaluation (masks off invalid depth pixels), unlike what Wonka et.al. do in their paper."
$ws.Cells.Item(37, 12).WrapText = $true
$ws.Cells.Item(37, 12).Value = "Pytorch version. Will run caffe version soon."

# ---- Row heights (wrapped notes make these rows taller, as in rows 3-5) ----
$ws.Rows.Item(37).RowHeight = 17
$ws.Rows.Item(38).RowHeight = 51

# ---- View state: scroll/selection ends up on the newly added note cell ----
$ws.Range("L38").Select()

Write-Output "done"
